$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from H1 (header style) into the new header cells I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header labels for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for rows 2-46: row, I value, J value
$data = @(
    @(2, 6, 6),
    @(3, 8, 9),
    @(4, 7, 7),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 10, 10),
    @(8, 7, 7),
    @(9, 5, 7),
    @(10, 8, 8),
    @(11, 6, 8),
    @(12, 5, 7),
    @(13, 5, 7),
    @(14, 1, 3),
    @(15, 1, 4),
    @(16, 4, 6),
    @(17, 6, 7),
    @(18, 10, 10),
    @(19, 6, 6),
    @(20, 8, 8),
    @(21, 6, 7),
    @(22, 7, 7),
    @(23, 6, 7),
    @(24, 5, 6),
    @(25, 5, 5),
    @(26, 6, 7),
    @(27, 5, 5),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 6, 6),
    @(31, 7, 8),
    @(32, 5, 5),
    @(33, 8, 8),
    @(34, 10, 10),
    @(35, 7, 7),
    @(36, 8, 9),
    @(37, 7, 7),
    @(38, 10, 10),
    @(39, 6, 7),
    @(40, 7, 7),
    @(41, 6, 7),
    @(42, 7, 7),
    @(43, 7, 8),
    @(44, 5, 7),
    @(45, 4, 6),
    @(46, 3, 4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
